$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = 0.20132366145238353
